# maj images _ ref  2024
# Fills in previously-blank "Origine expérience" / "Client" cells for four
# existing rows, flips two "(en cours)" placeholder years to their final
# literal year now that those missions wrapped up, and appends four brand
# new reference rows (43-46). Finishes by turning the data range into a
# proper AutoFilter (with its hidden _FilterDatabase name) now that the
# sheet has grown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFormat = '_-* #,##0\ "€"_-;\-* #,##0\ "€"_-;_-* "-"??\ "€"_-;_-@_-'

# --- Row 39: DDETS du Var (diagnostic jeunes sans abri) ---------------------
$ws.Range("A39").Value = "d-sidd"
$ws.Range("B39").Value = "DDETS du Var"

# --- Row 40: Logi-cité (logement social Gap Tallard Durance) ---------------
$ws.Range("A40").Value = "d-sidd"
$ws.Range("B40").Value = "Logi-cité"
$ws.Range("F40").Value = 2023

# --- Row 41: Terre d'Avance (PAM 74) ----------------------------------------
$ws.Range("A41").Value = "d-sidd"
$ws.Range("B41").Value = "Terre d'Avance"

# --- Row 42: CERC Ile de France (tableaux de bord BTP IDF) -----------------
$ws.Range("A42").Value = "d-sidd"
$ws.Range("B42").Value = "CERC Ile de France"
$ws.Range("F42").Value = 2023

# --- Row 43: Département Allier ---------------------------------------------
$ws.Range("A43").Value = "d-sidd"
$ws.Range("B43").Value = "Département Allier"
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = "Développement économique durable des territoires"
$ws.Range("E43").Value = "Cartographie interactive pour l’accompagnement au diagnostic stratégique préalable d’un pacte local des solidarités dans le département de l’Allier "
$ws.Range("F43").Value = 2023
$ws.Range("G43").Value = 3400
$ws.Range("G43").NumberFormat = $currencyFormat
$ws.Range("I43").Value = "Oui"

# --- Row 44: CERC Occitanie (Observatoire Immobilier d'Entreprise) ---------
$ws.Range("A44").Value = "d-sidd"
$ws.Range("B44").Value = "CERC Occitanie"
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = "Développement économique durable des territoires"
$ws.Range("E44").Value = "Observatoire de l' Immobilier d'Entreprise de Montpellier et son territoire Urbain"
$ws.Range("F44").Value = "2022, 2023, 2024"
$ws.Range("G44").Formula = "=6300*2"
$ws.Range("G44").NumberFormat = $currencyFormat
$ws.Range("I44").Value = "Oui"

# --- Row 45: CERC Occitanie (Observatoire Immobilier du Commerce) ----------
$ws.Range("A45").Value = "d-sidd"
$ws.Range("B45").Value = "CERC Occitanie"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = "Développement économique durable des territoires"
$ws.Range("E45").Value = "Observatoire de l'Immobilier du Commerce Montpellier et son territoire Urbain"
$ws.Range("F45").Value = 2022
$ws.Range("G45").Value = 6300
$ws.Range("G45").NumberFormat = $currencyFormat
$ws.Range("I45").Value = "Oui"

# --- Row 46: CEISS Consultants ----------------------------------------------
$ws.Range("A46").Value = "d-sidd"
$ws.Range("B46").Value = "CEISS Consultants"
$ws.Range("C46").Value = 1
$ws.Range("D46").Value = "Développement économique durable des territoires"
$ws.Range("E46").Value = "Données de cadrage du plan départemental d'action pour le logement et l'hébergement des personnes défavorisées des Alpes de Hautes Provence"
$ws.Range("F46").Value = 2022
$ws.Range("G46").Value = 3570
$ws.Range("G46").NumberFormat = $currencyFormat
$ws.Range("I46").Value = "Oui"
$ws.Range("J46").Value = "Oui"

# --- Turn the (now larger) data range into an AutoFilter --------------------
$ws.Range("A1:J46").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Feuil1!`$A`$1:`$J`$46")
$filterName.Visible = $false

# --- Restore selection near the newly-added rows -----------------------------
$ws.Range("K37").Select()
